$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column M (13th column), shifting M:AD to N:AE
$ws.Columns("M:M").Insert()

# Fill in the new values
$ws.Range("M6").Value = "1,25 timer"
$ws.Range("M7").Value = "1 time"
$ws.Range("L8").Value = "formel beskriv UC3"
$ws.Range("M8").Value = "1 time"
$ws.Range("L9").Value = "skrive rapport om datamodel og normalisering"
$ws.Range("M9").Value = "1 time"

Write-Host "done"
